$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension-relevant data: rows 2 to 23, columns A to I

# Row 2: Baahubali Crown of Blood
$ws.Range("A2").Value = 'Baahubali Crown of Blood'
$ws.Range("B2").Value = "'2024"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = '1080p'
$ws.Range("D2").Value = 'AVC'
$ws.Range("E2").Value = 'Hin + Kan + Mal + Tam + Tel'
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = "'01"
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value = "'06"
$ws.Range("H2").Style = "Normal"
$ws.Range("I2").Value = '590.18 MB'

# Row 3: Baahubali Crown of Blood
$ws.Range("A3").Value = 'Baahubali Crown of Blood'
$ws.Range("B3").Value = "'2024"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = '1080p'
$ws.Range("D3").Value = 'AVC'
$ws.Range("E3").Value = 'Hin + Kan + Mal + Tam + Tel'
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = "'01"
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").Value = "'04"
$ws.Range("H3").Style = "Normal"
$ws.Range("I3").Value = '603.9 MB'

# Row 4: The Flash
$ws.Range("A4").Value = 'The Flash'
$ws.Range("B4").Value = "'2014"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = 'BluRay'
$ws.Range("D4").Value = 'x264'
$ws.Range("E4").Value = 'Eng + Hin + Tam'
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = "'01"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value = "'05"
$ws.Range("H4").Style = "Normal"
$ws.Range("I4").Value = '613.37 MB'

# Row 5: The Flash
$ws.Range("A5").Value = 'The Flash'
$ws.Range("B5").Value = "'2014"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = 'BluRay'
$ws.Range("D5").Value = 'x264'
$ws.Range("E5").Value = 'Eng + Hin + Tam'
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = "'01"
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value = "'04"
$ws.Range("H5").Style = "Normal"
$ws.Range("I5").Value = '612.95 MB'

# Row 6: The Flash
$ws.Range("A6").Value = 'The Flash'
$ws.Range("B6").Value = "'2014"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = 'BluRay'
$ws.Range("D6").Value = 'x264'
$ws.Range("E6").Value = 'Eng + Hin + Tam'
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = "'01"
$ws.Range("G6").Style = "Normal"
$ws.Range("H6").Value = "'01"
$ws.Range("H6").Style = "Normal"
$ws.Range("I6").Value = '641.95 MB'

# Row 7: The Flash
$ws.Range("A7").Value = 'The Flash'
$ws.Range("B7").Value = "'2014"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = 'BluRay'
$ws.Range("D7").Value = 'x264'
$ws.Range("E7").Value = 'Eng + Hin + Tam'
$ws.Range("F7").Value = ""
$ws.Range("G7").Value = "'01"
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").Value = "'02"
$ws.Range("H7").Style = "Normal"
$ws.Range("I7").Value = '613.14 MB'

# Row 8: The Flash
$ws.Range("A8").Value = 'The Flash'
$ws.Range("B8").Value = "'2014"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = 'BluRay'
$ws.Range("D8").Value = 'x264'
$ws.Range("E8").Value = 'Eng + Hin + Tam'
$ws.Range("F8").Value = ""
$ws.Range("G8").Value = "'01"
$ws.Range("G8").Style = "Normal"
$ws.Range("H8").Value = "'03"
$ws.Range("H8").Style = "Normal"
$ws.Range("I8").Value = '612.69 MB'

# Row 9: Baahubali Crown of Blood
$ws.Range("A9").Value = 'Baahubali Crown of Blood'
$ws.Range("B9").Value = "'2024"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = '1080p'
$ws.Range("D9").Value = 'AVC'
$ws.Range("E9").Value = 'Hin + Kan + Mal + Tam + Tel'
$ws.Range("F9").Value = ""
$ws.Range("G9").Value = "'01"
$ws.Range("G9").Style = "Normal"
$ws.Range("H9").Value = "'03"
$ws.Range("H9").Style = "Normal"
$ws.Range("I9").Value = '601.64 MB'

# Row 10: Deadpool & Wolverine
$ws.Range("A10").Value = 'Deadpool & Wolverine'
$ws.Range("B10").Value = "'2024"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = '4K'
$ws.Range("D10").Value = 'SDR'
$ws.Range("E10").Value = 'English'
$ws.Range("F10").Value = '768Kbps & AAC & DD+5.1'
$ws.Range("G10").Value = ""
$ws.Range("H10").Value = ""
$ws.Range("I10").Value = '22.53 GB'

# Row 11: Deadpool & Wolverine
$ws.Range("A11").Value = 'Deadpool & Wolverine'
$ws.Range("B11").Value = "'2024"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = 'BluRay'
$ws.Range("D11").Value = 'x264'
$ws.Range("E11").Value = 'Eng + Hin + Tam + Tel'
$ws.Range("F11").Value = '192Kbps & DD+5.1'
$ws.Range("G11").Value = ""
$ws.Range("H11").Value = ""
$ws.Range("I11").Value = '1.66 GB'

# Row 12: Meiyazhagan
$ws.Range("A12").Value = 'Meiyazhagan'
$ws.Range("B12").Value = "'2024"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = 'HQ HDRip'
$ws.Range("D12").Value = 'x264'
$ws.Range("E12").Value = 'Tamil + Tamil UNCUT - HQ HDRip - x264 - AAC - 700MB - ESub'
$ws.Range("F12").Value = 'AAC'
$ws.Range("G12").Value = ""
$ws.Range("H12").Value = ""
$ws.Range("I12").Value = '739.0 MB'

# Row 13: Meiyazhagan
$ws.Range("A13").Value = 'Meiyazhagan'
$ws.Range("B13").Value = "'2024"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = 'HQ HDRip'
$ws.Range("D13").Value = 'x264'
$ws.Range("E13").Value = 'Tamil + Tamil - HQ HDRip - x264 - AAC - 700MB - ESub'
$ws.Range("F13").Value = 'AAC'
$ws.Range("G13").Value = ""
$ws.Range("H13").Value = ""
$ws.Range("I13").Value = '739.14 MB'

# Row 14: Pushpa 2 The Rule
$ws.Range("A14").Value = 'Pushpa 2 The Rule'
$ws.Range("B14").Value = "'2024"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = '720p'
$ws.Range("D14").Value = 'x264'
$ws.Range("E14").Value = 'Kannada + Kannada - 720p HDRip - x264 - AAC - 1.4GB - Original Audio'
$ws.Range("F14").Value = 'AAC'
$ws.Range("G14").Value = ""
$ws.Range("H14").Value = ""
$ws.Range("I14").Value = '1.42 GB'

# Row 15: Pushpa 2 The Rule
$ws.Range("A15").Value = 'Pushpa 2 The Rule'
$ws.Range("B15").Value = "'2024"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = '1080p'
$ws.Range("D15").Value = 'AVC'
$ws.Range("E15").Value = 'Tamil + Tamil - 1080p HD AVC UNTOUCHED - x264 - AAC - 4.8GB'
$ws.Range("F15").Value = 'AAC'
$ws.Range("G15").Value = ""
$ws.Range("H15").Value = ""
$ws.Range("I15").Value = '4.84 GB'

# Row 16: Miss You
$ws.Range("A16").Value = 'Miss You'
$ws.Range("B16").Value = "'2024"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = 'HQ HDRip'
$ws.Range("D16").Value = 'x265'
$ws.Range("E16").Value = 'Tamil'
$ws.Range("F16").Value = 'AAC'
$ws.Range("G16").Value = ""
$ws.Range("H16").Value = ""
$ws.Range("I16").Value = '909.73 MB'

# Row 17: Miss You
$ws.Range("A17").Value = 'Miss You'
$ws.Range("B17").Value = "'2024"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = '720p'
$ws.Range("D17").Value = 'x264'
$ws.Range("E17").Value = 'AAC + Tamil + Tamil - 720p HQ HDRip - x264 - [DD5.1(192Kbps)'
$ws.Range("F17").Value = 'AAC & 192Kbps'
$ws.Range("G17").Value = ""
$ws.Range("H17").Value = ""
$ws.Range("I17").Value = '1.14 GB'

# Row 18: Singham Again
$ws.Range("A18").Value = 'Singham Again'
$ws.Range("B18").Value = "'2024"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = 'HQ HDRip'
$ws.Range("D18").Value = 'x264'
$ws.Range("E18").Value = 'Hindi + Hindi - HQ HDRip - x264 - AAC - 800MB'
$ws.Range("F18").Value = 'AAC'
$ws.Range("G18").Value = ""
$ws.Range("H18").Value = ""
$ws.Range("I18").Value = '837.42 MB'

# Row 19: Singham Again
$ws.Range("A19").Value = 'Singham Again'
$ws.Range("B19").Value = "'2024"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = '720p'
$ws.Range("D19").Value = 'HEVC'
$ws.Range("E19").Value = 'AAC + Tamil + Tamil - 720p HQ HDRip HEVC - x265 - [DD5.1(192Kbps)'
$ws.Range("F19").Value = 'AAC & 192Kbps'
$ws.Range("G19").Value = ""
$ws.Range("H19").Value = ""
$ws.Range("I19").Value = '984.75 MB'

# Row 20: Barroz
$ws.Range("A20").Value = 'Barroz'
$ws.Range("B20").Value = "'2024"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = '720p'
$ws.Range("D20").Value = 'HEVC'
$ws.Range("E20").Value = 'Tamil'
$ws.Range("F20").Value = 'AAC'
$ws.Range("G20").Value = ""
$ws.Range("H20").Value = ""
$ws.Range("I20").Value = '972.66 MB'

# Row 21: Barroz
$ws.Range("A21").Value = 'Barroz'
$ws.Range("B21").Value = "'2024"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = 'HQ HDRip'
$ws.Range("D21").Value = 'HEVC'
$ws.Range("E21").Value = 'Tamil'
$ws.Range("F21").Value = 'AAC & 192Kbps & DD+5.1'
$ws.Range("G21").Value = ""
$ws.Range("H21").Value = ""
$ws.Range("I21").Value = '1.1 GB'

# Row 22: Game Changer
$ws.Range("A22").Value = 'Game Changer'
$ws.Range("B22").Value = "'2025"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = '720p'
$ws.Range("D22").Value = 'HEVC'
$ws.Range("E22").Value = 'AAC + Tamil + Tamil - 720p HQ HDRip HEVC - x265 - [DD5.1(192Kbps)'
$ws.Range("F22").Value = 'AAC & 192Kbps'
$ws.Range("G22").Value = ""
$ws.Range("H22").Value = ""
$ws.Range("I22").Value = '1.26 GB'

# Row 23: Game Changer
$ws.Range("A23").Value = 'Game Changer'
$ws.Range("B23").Value = "'2025"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = 'HQ HDRip'
$ws.Range("D23").Value = 'HEVC'
$ws.Range("E23").Value = 'Tamil'
$ws.Range("F23").Value = 'AAC'
$ws.Range("G23").Value = ""
$ws.Range("H23").Value = ""
$ws.Range("I23").Value = '957.45 MB'
